$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45934
$ws.Range("B2").Value = 110.94
$ws.Range("C2").Value = 102.92
$ws.Range("D2").Value = 100.53
$ws.Range("E2").Value = 95.75
$ws.Range("F2").Value = 93.01000000000001
$ws.Range("G2").Value = 99.63
$ws.Range("H2").Value = 94.38
$ws.Range("I2").Value = 100.3
$ws.Range("J2").Value = 92.34999999999999
$ws.Range("K2").Value = 40.37
$ws.Range("L2").Value = 6.2
$ws.Range("M2").Value = 0.18
$ws.Range("N2").Value = 0.08
$ws.Range("O2").Value = 0.21
$ws.Range("P2").Value = 0.07000000000000001
$ws.Range("Q2").Value = -0.01
$ws.Range("R2").Value = -0.01
$ws.Range("S2").Value = 0.15
$ws.Range("T2").Value = 10.81
$ws.Range("U2").Value = 65.61
$ws.Range("V2").Value = 74.34999999999999
$ws.Range("W2").Value = 71.95
$ws.Range("X2").Value = 68.26000000000001
$ws.Range("Y2").Value = 47.03
$ws.Range("Z2").Value = 53.13
$ws.Range("AA2").Value = "0h-4h"
$ws.Range("AB2").Value = 102.54
$ws.Range("AC2").Value = "0h-2h"
$ws.Range("AD2").Value = 106.93
$ws.Range("AE2").Value = "2h-4h"
$ws.Range("AF2").Value = 98.14
$ws.Range("AG2").Value = "9h-23h"
